$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the Targum columns
$ws.Range("J1").Value = "Onkelos"
$ws.Range("K1").Value = "Jonathan"
$ws.Range("J1:K1").Style = $ws.Range("A1").Style

# New Targum text, written column-major (J2, J3, K2, K3) so new shared
# strings are created in the same order as the target workbook.
$ws.Range("J2").Value = "And there has not ever arisen a prophet within Yisroel like Moshe, whom Adonoy knew [<b>appeared to</b>] face-to-face."
$ws.Range("J3").Value = "“Go, gather the elders of Yisrael, and say to them, ‘Adonoy, the God of your fathers appeared [<b>became revealed</b>] to me—the God of Avraham, Yitzchok and Yaakov—saying, “I have indeed been mindful of you, regarding that which is being done to you in Egypt."
$ws.Range("K2").Value = "But no prophet hath again risen in Israel like unto Mosheh, because the Word of the Lord had known him to speak with him word for word,"
$ws.Range("K3").Value = "Go, and assemble the elders of Israel, and say to them, The Lord God of your fathers hath appeared unto me, the God of Abraham, Izhak, and Jakob, saying, Remembering, I have remembered you, and the injury that is done you in Mizraim;"
